$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.169.88"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "1.824.72"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'241.54"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").Value = "'0.6194"
$ws.Range("E6").Value = "  -0.97%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.07350"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").Value = "'0.2897"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("D11").Value = "'0.07669"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "1.821.74"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").Value = "'0.6625"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "'82.26"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").Value = "'0.000008950"
$ws.Range("E16").Value = "  -4.52%  "

$ws.Range("D17").Value = "'5.835"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("D18").Value = "29.162.61"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").Value = "2.070.02"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'238.24"
$ws.Range("E20").Value = "  +6.90%  "

$ws.Range("D21").Value = "'12.42"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "'7.199"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "'158.01"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").Value = "'0.1419"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").Value = "'8.454"
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("D28").Value = "'17.64"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("D29").Value = "'1.482"
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").Value = "'0.05587"
$ws.Range("E30").Value = "  -4.45%  "

$ws.Range("D31").Value = "'4.091"
$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").Value = "'4.094"
$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("D33").Value = "'1.205"
$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("D34").Value = "'1.833"
$ws.Range("E34").Value = "  +0.25%  "

$ws.Range("D35").Value = "'0.7351"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").Value = "'1.133"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D37").Value = "'2.628"
$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").Value = "'2.840"
$ws.Range("E38").Value = "  +2.81%  "

$ws.Range("D39").Value = "1.220.10"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").Value = "'0.01763"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.312"
$ws.Range("E41").Value = "  -2.73%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9157"
$ws.Range("E42").Value = "  +2.82%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'101.59"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").Value = "1.973.75"
$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("D46").Value = "'64.77"
$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("D47").Value = "'0.5080"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("E48").Value = "  -6.29%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.104"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4011"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").Value = "'0.05762"
$ws.Range("E51").Value = "  -1.05%  "
